# fix: minor changes in slides for the coming semesters
$p = $ppt.ActivePresentation

# --- Slide 12: rename referenced notebook file ---
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(4)
$tr12 = $shp12.TextFrame.TextRange
$found = $tr12.Find("overfitting.ipynb")
$found.Text = "overfitting_experimentation.ipynb"

# --- Slide 20: remove the "trophy" announcement paragraph from the
#     homework body text ---
$s20 = $p.Slides.Item(20)
$shp20 = $s20.Shapes.Item(4)
$tr20 = $shp20.TextFrame.TextRange

$trophyPara = $tr20.Paragraphs(7, 1)
$trophyPara.Delete()
